$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Book1")

# --- Block 3 (rows 14-18): task wraps to next day, add "Slöläs" entry ---
$ws.Range("A17").Value = "Slöläs"
$ws.Range("B17").Value = 2/24
$ws.Range("B18").Value = 2.5/24
$ws.Range("B18").NumberFormat = "h:mm:ss"

# --- Block 1 (rows 1-5): add a third pair of columns (E/F) for "Jobb 09" ---
# header
$ws.Range("E1").Value = "Jobb 09"

# task "Annat" started at 22:00 (previous day) and wrapped to 10:00 next day
$ws.Range("E2").Value = "Annat"
$ws.Range("F2").Value = 10/24
$ws.Range("F2").NumberFormat = "h:mm:ss"

# task "Kaffe" at 20:00
$ws.Range("E3").Value = "Kaffe"
$ws.Range("F3").Value = 20/24
$ws.Range("F3").NumberFormat = "h:mm:ss"

# task "Sova" at 01:00
$ws.Range("E4").Value = "Sova"
$ws.Range("F4").Value = 1/24
$ws.Range("F4").NumberFormat = "h:mm:ss"

# Move selection to F5 (no explicit scroll position)
$null = $ws.Range("F5").Select()
